# Applies "Improves sentence structures for report chapters" edits.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

$rsq = [char]0x2019   # right single quotation mark '
$lsq = [char]0x2018   # left single quotation mark '

# --- Paragraph: "User Experience (UX) Design plays a vital role ..." ---

Replace-Text "and your underlying database" "and when underlying database"
Replace-Text "better UI and also usable? Its Design cycle as seen in the illustrates" "better UI and also usable? Its Design cycle, as seen in the illustrates"
Replace-Text "where we gather requirements, then we make a prototype based on the requirements and then tested, then again we gather new requirements based on the evaluation." "where we gather requirements. Next, we make a prototype based on the requirements and then tested, then again, we gather new requirements based on the evaluation."

# --- "We can see the applicability ..." ---
Replace-Text "We can see the applicability of the User Experience Design in detail at ." "We can see the applicability of the User Experience Design in detail at."

# --- "The target users for evaluation ..." ---
Replace-Text "Computer Science. This ensures that the evaluation process" "Computer Science. This users qualification ensures that the evaluation process"

# --- "Where N represents ... L represents ..." ---
Replace-Text "that after the number of users is five, then the usability problems" "that after the number of users is five, the usability problems"

# --- "As the order of prototypes ..." ---
Replace-Text "as they tend to learn. Therefore, the order is changed for different segments" "as they tend to learn. Therefore, we change the order for different segments"

# --- "In a cognitive walkthrough, ..." ---
Replace-Text "Blackmon, Polson, et al. in their paper  mentions four questions" "Blackmon, Polson, et al. in their paper mentions four questions"

# --- Cognitive walkthrough 4 questions (list items) ---
Replace-Text "Will the user try and achieve the right outcome?" "Will the user be able to try and attain the right conclusion?"
Replace-Text "Will the user notice that the correct action is available to them?" "Will the user be able to notice that the right action is presented?"
Replace-Text "Will the user associate the correct action with the outcome they expect to achieve?" "Will the user be able to associate the right action with the outcome they expect to accomplish?"
Replace-Text "If the correct action is performed; will the user see that progress is being made towards their intended outcome?" "If the right action is made; will the user be able to see that progress is being made towards their intended conclusion?"

# --- Move the _GoBack bookmark from list item 1 to the end of list item 4 ---
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()
$q = $d.Content
$qFound = $q.Find.Execute("If the right action is made; will the user be able to see that progress is being made towards their intended conclusion?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($q.End, $q.End)
$bmAdded = $d.Bookmarks.Add("_GoBack", $bmRange)

# --- "These questions are also quite applicable ..." paragraph ---
Replace-Text "These questions are also quite applicable in our context. Thereby, these questions are assessed for each step, which is predetermined with designed elements on the user interface in order to solve the research question been tackled. So, the steps vary for each design. This approach gives qualitative feedback from a user as they are a mostly open-ended scenario to discuss primarily, for questions which are answered as" "These questions are also quite applicable in our context, and so, we assess these questions for each step. The designed elements on the user interface predetermines in order to solve a research question. So, the steps vary for each design. This approach gives qualitative feedback from a user as they are a mostly open-ended scenario to discuss primarily, for questions which users answers as"

# --- "Overall, cognitive walkthrough ..." ---
Replace-Text "Further, when over two best solution ideas are needed to be evaluated against each other" "Further, when we need to evaluate two best solution ideas against each other"

# --- Likert scale "Uni Polar Likert Scale shows one attribute ..." paragraph (full rewrite) ---
$oldUni = "Uni Polar Likert Scale shows one attribute let" + $rsq + "s say strongly agree on end and I do not agree on the other end whereas for Bi-Polar Likert Scale it shows strongly agree one side and strongly disagree on the other side and also there is a label in middle which usually states " + $lsq + "neither" + $rsq + ". In this thesis work, the Uni Polar Likert Scale is used with labels 0 and 10 on either end of a scale showing " + $lsq + "not at all usable" + $rsq + " to " + $lsq + "highly usable" + $rsq + ". Not at all usable could be understood as the worst design or solution idea for an user scenario. "
$newUni = "Uni Polar Likert Scale shows one attribute let us say " + $lsq + "strongly agree" + $rsq + " on one end and " + $lsq + "I do not agree" + $rsq + " on the other end. Incase of Bi-Polar Likert Scale it shows " + $lsq + "strongly agree" + $rsq + " on one side and " + $lsq + "strongly disagree" + $rsq + " on the other side, and also there is a label in the middle which usually states " + $lsq + "neither" + $rsq + ". In this thesis work, we use the " + $lsq + "Uni Polar Likert Scale" + $rsq + " with labels 0 and 10 on either end of a scale showing " + $lsq + "not at all usable" + $rsq + " to " + $lsq + "highly usable" + $rsq + ". Not at all usable could be understood as the worst design or solution idea for an user scenario. "
Replace-Text $oldUni $newUni

Write-Output "done-part1"
